$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A39").Value = "HFTRE4"
$ws.Range("B39").Value = "Cuchilla de limpieza"
$ws.Range("C39").Value = "TS3525"
$ws.Range("D39").Value = 0
$ws.Range("E39").Value = 100000
$ws.Range("F39").Value = 1
$ws.Range("G39").Value = 0
$ws.Range("H39").Formula = "=(E39-D39)*G39"
$ws.Range("I39").Formula = "=D39*F39"
$ws.Range("J39").Value = 0
